$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.761.99"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "2.350.28"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("E5").Value = "  +3.32%  "
$ws.Range("E6").Value = "  +0.93%  "
$ws.Range("E7").Value = "  +11.38%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  +19.00%  "
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("E11").Value = "  +3.62%  "
$ws.Range("E12").Value = "  +1.80%  "
$ws.Range("D13").Value = "2.705.41"
$ws.Range("E13").Value = "  +0.80%  "
$ws.Range("E14").Value = "  +7.71%  "
$ws.Range("E15").Value = "  +7.80%  "
$ws.Range("E16").Value = "  +4.23%  "
$ws.Range("D17").Value = "2.381.84"
$ws.Range("E17").Value = "  +1.88%  "
$ws.Range("D18").Value = "43.721.63"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("E19").Value = "  +3.34%  "
$ws.Range("E20").Value = "  +4.26%  "
$ws.Range("E21").Value = "  +2.54%  "
$ws.Range("E22").Value = "  +1.55%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  -1.36%  "
$ws.Range("E25").Value = "  +1.81%  "
$ws.Range("E26").Value = "  +6.46%  "
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("E29").Value = "  +9.32%  "
$ws.Range("E30").Value = "  -1.72%  "
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("E32").Value = "  +5.15%  "
$ws.Range("E33").Value = "  +2.95%  "
$ws.Range("E34").Value = "  +4.35%  "
$ws.Range("E35").Value = "  +3.79%  "
$ws.Range("E36").Value = "  +6.11%  "
$ws.Range("E37").Value = "  -1.84%  "
$ws.Range("E38").Value = "  -2.60%  "
$ws.Range("E39").Value = "  +6.12%  "
$ws.Range("E40").Value = "  +6.14%  "
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("E42").Value = "  -1.61%  "
$ws.Range("E43").Value = "  +2.22%  "
$ws.Range("E44").Value = "  -1.59%  "
$ws.Range("E45").Value = "  +1.67%  "
$ws.Range("E46").Value = "  +13.33%  "
$ws.Range("E47").Value = "  +2.37%  "
$ws.Range("E48").Value = "  -1.89%  "
$ws.Range("D49").Value = "1.432.02"
$ws.Range("E49").Value = "  -0.99%  "
$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("E50").Value = "  +1.63%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.577.85"
$ws.Range("E51").Value = "  +0.84%  "

# Numeric-looking price strings must stay as text; force text format, set, then reset style
$numericTextCells = @("D6","D7","D9","D10","D11","D14","D15","D20","D22","D24","D25","D26","D27","D28","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D42","D43","D46","D48","D50")
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range("D6").Value = "235.25"
$ws.Range("D7").Value = "73.33"
$ws.Range("D9").Value = "0.548"
$ws.Range("D10").Value = "0.0986"
$ws.Range("D11").Value = "27.95"
$ws.Range("D14").Value = "16.67"
$ws.Range("D15").Value = "6.69"
$ws.Range("D20").Value = "77.26"
$ws.Range("D22").Value = "253.18"
$ws.Range("D24").Value = "3.75"
$ws.Range("D25").Value = "2.48"
$ws.Range("D26").Value = "10.58"
$ws.Range("D27").Value = "2.28"
$ws.Range("D28").Value = "22.32"
$ws.Range("D30").Value = "172.12"
$ws.Range("D31").Value = "0.129"
$ws.Range("D32").Value = "0.132"
$ws.Range("D33").Value = "5.16"
$ws.Range("D34").Value = "0.0718"
$ws.Range("D35").Value = "5.17"
$ws.Range("D36").Value = "3.83"
$ws.Range("D37").Value = "2.42"
$ws.Range("D38").Value = "6.38"
$ws.Range("D39").Value = "0.0269"
$ws.Range("D40").Value = "19.35"
$ws.Range("D42").Value = "8.89"
$ws.Range("D43").Value = "0.0977"
$ws.Range("D46").Value = "0.184"
$ws.Range("D48").Value = "97.32"
$ws.Range("D50").Value = "2.78"
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).Style = "Normal"
}
